$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.4866
$ws.Range("A3").Value = -22.1279
$ws.Range("E3").Value = 16.34229999999999
$ws.Range("E12").Value = 18.32170000000002
$ws.Range("A14").Value = -21.8088
$ws.Range("A21").Value = -20.12289999999998
$ws.Range("A23").Value = -20.42659999999997
$ws.Range("E24").Value = 16.6845
$ws.Range("A25").Value = -21.92219999999999
$ws.Range("C25").Value = -12.7445
$ws.Range("E25").Value = 16.90300000000001
$ws.Range("A26").Value = -21.10709999999997
$ws.Range("C27").Value = -12.6782
$ws.Range("A29").Value = -20.98529999999997
$ws.Range("C31").Value = -13.13710000000001
$ws.Range("C39").Value = -12.84790000000001
$ws.Range("C48").Value = -11.6778
$ws.Range("E50").Value = 16.2838
$ws.Range("C51").Value = -11.8844
$ws.Range("C52").Value = -11.5092
$ws.Range("A53").Value = -21.98499999999999
$ws.Range("E53").Value = 16.72670000000002
$ws.Range("C55").Value = -13.76489999999999
$ws.Range("C56").Value = -12.3343
$ws.Range("A57").Value = -22.0871
$ws.Range("C57").Value = -13.057
$ws.Range("E57").Value = 16.47539999999999
$ws.Range("A59").Value = -22.2221
$ws.Range("E61").Value = 16.5998
$ws.Range("E63").Value = 18.24560000000002
$ws.Range("A69").Value = -21.555
$ws.Range("E70").Value = 18.40880000000002
$ws.Range("C73").Value = -13.0158
$ws.Range("A79").Value = -20.33150000000001
$ws.Range("A83").Value = -21.8599
$ws.Range("E86").Value = 16.4747
$ws.Range("C89").Value = -10.0055
$ws.Range("C90").Value = -11.79510000000001
$ws.Range("A91").Value = -21.4996
$ws.Range("C92").Value = -10.38169999999999
$ws.Range("A93").Value = -20.86199999999998
$ws.Range("E98").Value = 15.9275
$ws.Range("E100").Value = 16.36620000000001
$ws.Range("E102").Value = 16.78009999999999
